# Update the answer text in each table cell of the "three-digit divided by
# one-digit" worksheet to the newly generated set of division problems.
# Each old answer string is unique within the document, so a plain
# whole-document Find/Replace (MatchCase, MatchWholeWord) safely targets
# exactly one run per call.

$d = $word.ActiveDocument

$replacements = @(
    @("230÷9=25, 5", "296÷5=59, 1"),
    @("316÷2=158, 0", "868÷2=434, 0"),
    @("266÷9=29, 5", "501÷6=83, 3"),
    @("856÷7=122, 2", "361÷7=51, 4"),
    @("233÷8=29, 1", "159÷9=17, 6"),
    @("910÷3=303, 1", "585÷9=65, 0"),
    @("558÷9=62, 0", "338÷2=169, 0"),
    @("421÷8=52, 5", "815÷7=116, 3"),
    @("194÷3=64, 2", "590÷7=84, 2"),
    @("782÷6=130, 2", "393÷7=56, 1"),
    @("222÷2=111, 0", "992÷5=198, 2"),
    @("582÷6=97, 0", "750÷6=125, 0"),
    @("522÷3=174, 0", "890÷8=111, 2"),
    @("108÷9=12, 0", "106÷7=15, 1"),
    @("280÷9=31, 1", "870÷4=217, 2"),
    @("759÷8=94, 7", "940÷6=156, 4"),
    @("334÷5=66, 4", "866÷7=123, 5"),
    @("871÷6=145, 1", "653÷6=108, 5"),
    @("133÷7=19, 0", "477÷9=53, 0"),
    @("964÷4=241, 0", "883÷3=294, 1"),
    @("123÷9=13, 6", "735÷6=122, 3"),
    @("483÷3=161, 0", "552÷9=61, 3"),
    @("391÷8=48, 7", "914÷7=130, 4"),
    @("346÷4=86, 2", "577÷4=144, 1"),
    @("629÷7=89, 6", "692÷8=86, 4")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: replacement not applied, old text not found: $old"
    }
}

Write-Output "Done."
